$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2960896.8
$ws.Range("J17").Value = 3050563.2
$ws.Range("L17").Value = 9151689.600000001
$ws.Range("N17").Value = -9152025.600000001
$ws.Range("H43").Value = 3232.1428
$ws.Range("I43").Value = 672.5
$ws.Range("J43").Value = 4807.3076
$ws.Range("K43").Value = 672.5
$ws.Range("L43").Value = 4807.3076
$ws.Range("M43").Value = -603.5
$ws.Range("N43").Value = -4945.3076
$ws.Range("H99").Value = 268.75
$ws.Range("I99").Value = 225
$ws.Range("K99").Value = 675
$ws.Range("M99").Value = 823
$ws.Range("H101").Value = 285
$ws.Range("I101").Value = 282
$ws.Range("J101").Value = 300
$ws.Range("K101").Value = 846
$ws.Range("L101").Value = 900
$ws.Range("M101").Value = 776
$ws.Range("N101").Value = -4144
$ws.Range("H113").Value = 55559140
$ws.Range("I113").Value = 125001660
$ws.Range("J113").Value = 5118.3
$ws.Range("K113").Value = 125001660
$ws.Range("L113").Value = 5118.3
$ws.Range("M113").Value = -124998406
$ws.Range("N113").Value = -11626.3
$ws.Range("H125").Value = 425
$ws.Range("I125").Value = 369.85715
$ws.Range("J125").Value = 502.2
$ws.Range("K125").Value = 3328.71435
$ws.Range("L125").Value = 4519.8
$ws.Range("M125").Value = -868.7143499999997
$ws.Range("N125").Value = -9439.799999999999
$ws.Range("H129").Value = 250934.2
$ws.Range("I129").Value = 298.5
$ws.Range("K129").Value = 895.5
$ws.Range("M129").Value = 4104.5
$ws.Range("H137").Value = 99244.61
$ws.Range("I137").Value = 155306.89
$ws.Range("J137").Value = 2070
$ws.Range("K137").Value = 465920.67
$ws.Range("L137").Value = 6210
$ws.Range("M137").Value = -463370.67
$ws.Range("N137").Value = -11310
$ws.Range("H138").Value = 3234.738
$ws.Range("I138").Value = 1995.8096
$ws.Range("J138").Value = 4473.6665
$ws.Range("K138").Value = 5987.4288
$ws.Range("L138").Value = 13420.9995
$ws.Range("M138").Value = -847.4287999999997
$ws.Range("N138").Value = -23700.9995

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9293.671
$ws.Range("I32").Value = 6913.027
$ws.Range("J32").Value = 21877.072
$ws.Range("K32").Value = 6913.027
$ws.Range("L32").Value = 21877.072
$ws.Range("M32").Value = -6626.027
$ws.Range("N32").Value = -22451.072
$ws.Range("H61").Value = 9527055
$ws.Range("I61").Value = 11114027
$ws.Range("K61").Value = 11114027
$ws.Range("M61").Value = -11113815
$ws.Range("H63").Value = 2841956.8
$ws.Range("I63").Value = 1147
$ws.Range("K63").Value = 1147
$ws.Range("M63").Value = -461
$ws.Range("H66").Value = 2841956.8
$ws.Range("I66").Value = 1147
$ws.Range("K66").Value = 5735
$ws.Range("M66").Value = -2303
$ws.Range("H74").Value = 23810854
$ws.Range("I74").Value = 35714904
$ws.Range("J74").Value = 2747.8572
$ws.Range("K74").Value = 35714904
$ws.Range("L74").Value = 2747.8572
$ws.Range("M74").Value = -35714030
$ws.Range("N74").Value = -4495.8572
$ws.Range("H77").Value = 23810854
$ws.Range("I77").Value = 35714904
$ws.Range("J77").Value = 2747.8572
$ws.Range("K77").Value = 178574520
$ws.Range("L77").Value = 13739.286
$ws.Range("M77").Value = -178570152
$ws.Range("N77").Value = -22475.286
$ws.Range("H132").Value = 11377322
$ws.Range("I132").Value = 13890953
$ws.Range("J132").Value = 65980
$ws.Range("K132").Value = 41672859
$ws.Range("L132").Value = 197940
$ws.Range("M132").Value = -41670329
$ws.Range("N132").Value = -203000
$ws.Range("H136").Value = 9527055
$ws.Range("I136").Value = 11114027
$ws.Range("K136").Value = 33342081
$ws.Range("M136").Value = -33339531

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 776.8
$ws.Range("I16").Value = 785.3333
$ws.Range("J16").Value = 700
$ws.Range("K16").Value = 785.3333
$ws.Range("L16").Value = 700
$ws.Range("M16").Value = -498.3333
$ws.Range("N16").Value = -1274
$ws.Range("H52").Value = 23948.889
$ws.Range("I52").Value = 8000
$ws.Range("J52").Value = 25942.5
$ws.Range("K52").Value = 8000
$ws.Range("L52").Value = 25942.5
$ws.Range("M52").Value = -7706
$ws.Range("N52").Value = -26530.5
$ws.Range("H62").Value = 2928.9756
$ws.Range("I62").Value = 2712.6843
$ws.Range("K62").Value = 2712.6843
$ws.Range("M62").Value = -2088.6843
$ws.Range("H65").Value = 2928.9756
$ws.Range("I65").Value = 2712.6843
$ws.Range("K65").Value = 13563.4215
$ws.Range("M65").Value = -10443.4215
$ws.Range("H86").Value = 12259.454
$ws.Range("I86").Value = 1556
$ws.Range("J86").Value = 21179
$ws.Range("K86").Value = 1556
$ws.Range("L86").Value = 21179
$ws.Range("M86").Value = -433
$ws.Range("N86").Value = -23425
$ws.Range("H89").Value = 12259.454
$ws.Range("I89").Value = 1556
$ws.Range("J89").Value = 21179
$ws.Range("K89").Value = 7780
$ws.Range("L89").Value = 105895
$ws.Range("M89").Value = -2164
$ws.Range("N89").Value = -117127
$ws.Range("H113").Value = 776.8
$ws.Range("I113").Value = 785.3333
$ws.Range("J113").Value = 700
$ws.Range("K113").Value = 785.3333
$ws.Range("L113").Value = 700
$ws.Range("M113").Value = 1384.6667
$ws.Range("N113").Value = -5040
$ws.Range("H122").Value = 1125.4242
$ws.Range("I122").Value = 932.2692
$ws.Range("K122").Value = 2796.8076
$ws.Range("M122").Value = -346.8076000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 781.93335
$ws.Range("I34").Value = 383.33334
$ws.Range("J34").Value = 881.5833
$ws.Range("K34").Value = 1150.00002
$ws.Range("L34").Value = 2644.7499
$ws.Range("M34").Value = -1066.00002
$ws.Range("N34").Value = -2812.7499
$ws.Range("H39").Value = 3500
$ws.Range("J39").Value = 3500
$ws.Range("L39").Value = 10500
$ws.Range("N39").Value = -11088
$ws.Range("H55").Value = 4249.75
$ws.Range("J55").Value = 4249.75
$ws.Range("L55").Value = 12749.25
$ws.Range("N55").Value = -13103.25
$ws.Range("H69").Value = 1975
$ws.Range("J69").Value = 1833.3334
$ws.Range("L69").Value = 5500.0002
$ws.Range("N69").Value = -7122.0002
$ws.Range("H70").Value = 3305.6
$ws.Range("I70").Value = 1357.1428
$ws.Range("J70").Value = 5010.5
$ws.Range("K70").Value = 4071.4284
$ws.Range("L70").Value = 15031.5
$ws.Range("M70").Value = -3756.4284
$ws.Range("N70").Value = -15661.5
$ws.Range("H72").Value = 1975
$ws.Range("J72").Value = 1833.3334
$ws.Range("L72").Value = 16500.0006
$ws.Range("N72").Value = -24612.0006
$ws.Range("H73").Value = 3305.6
$ws.Range("I73").Value = 1357.1428
$ws.Range("J73").Value = 5010.5
$ws.Range("K73").Value = 4071.4284
$ws.Range("L73").Value = 15031.5
$ws.Range("M73").Value = -2979.4284
$ws.Range("N73").Value = -17215.5
$ws.Range("H74").Value = 9969.4
$ws.Range("J74").Value = 9969.4
$ws.Range("L74").Value = 29908.2
$ws.Range("N74").Value = -32030.2
$ws.Range("H77").Value = 9969.4
$ws.Range("J77").Value = 9969.4
$ws.Range("L77").Value = 89724.59999999999
$ws.Range("N77").Value = -100332.6
$ws.Range("H82").Value = 9500
$ws.Range("I82").Value = 9000
$ws.Range("J82").Value = 10000
$ws.Range("K82").Value = 27000
$ws.Range("L82").Value = 30000
$ws.Range("M82").Value = -26594
$ws.Range("N82").Value = -30812
$ws.Range("H85").Value = 9500
$ws.Range("I85").Value = 9000
$ws.Range("J85").Value = 10000
$ws.Range("K85").Value = 27000
$ws.Range("L85").Value = 30000
$ws.Range("M85").Value = -25596
$ws.Range("N85").Value = -32808
$ws.Range("H113").Value = 774.2727
$ws.Range("J113").Value = 867.7692
$ws.Range("L113").Value = 2603.3076
$ws.Range("N113").Value = -6943.3076
$ws.Range("H122").Value = 1083.8182
$ws.Range("J122").Value = 1083.8182
$ws.Range("L122").Value = 9754.363799999999
$ws.Range("N122").Value = -14654.3638
$ws.Range("H131").Value = 699.33
$ws.Range("J131").Value = 741.86365
$ws.Range("L131").Value = 2225.59095
$ws.Range("N131").Value = -12305.59095

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4908393.5
$ws.Range("I132").Value = 6688960
$ws.Range("J132").Value = 75428.14
$ws.Range("K132").Value = 20066880
$ws.Range("L132").Value = 226284.42
$ws.Range("M132").Value = -20064350
$ws.Range("N132").Value = -231344.42
$ws.Range("H135").Value = 37207.8
$ws.Range("J135").Value = 37207.8
$ws.Range("L135").Value = 37207.8
$ws.Range("N135").Value = -47347.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 857.9231
$ws.Range("I16").Value = 856.375
$ws.Range("K16").Value = 856.375
$ws.Range("M16").Value = -686.375
$ws.Range("H93").Value = 1828.75
$ws.Range("I93").Value = 1866.1538
$ws.Range("K93").Value = 1866.1538
$ws.Range("M93").Value = -618.1538

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 186.27272
$ws.Range("I100").Value = 164.14285
$ws.Range("K100").Value = 328.2857
$ws.Range("M100").Value = 212.7143
$ws.Range("H107").Value = 233.46666
$ws.Range("I107").Value = 169.38461
$ws.Range("K107").Value = 508.15383
$ws.Range("M107").Value = 1411.84617
$ws.Range("H136").Value = 28574578
$ws.Range("I136").Value = 37038644
$ws.Range("K136").Value = 111115932
$ws.Range("M136").Value = -111113382
